# Insert 3 new weekly report rows (Hass - Especial/Primera/Segunda, fecha 2023-04-05)
# immediately above the existing row 615, shifting the old rows 615:664 down to 618:667.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 615 (this pushes rows 615:664 -> 618:667)
$ws.Rows("615:617").Insert()

# Shared/constant values for this data block (same across the whole sheet)
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$unidad    = "$/kilo (en caja de 17 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 1

# New rows data: row, fecha(serial), variedad, calidad, volumen, precioMin, precioMax, precioProm
$newRows = @(
    @{ Row = 615; Fecha = 45021; Variedad = "Hass"; Calidad = "Especial"; Volumen = 240; PMin = 3800; PMax = 3900; PProm = 3850 },
    @{ Row = 616; Fecha = 45021; Variedad = "Hass"; Calidad = "Primera";  Volumen = 300; PMin = 3500; PMax = 3600; PProm = 3550 },
    @{ Row = 617; Fecha = 45021; Variedad = "Hass"; Calidad = "Segunda";  Volumen = 200; PMin = 3200; PMax = 3300; PProm = 3250 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PProm
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
